# The document's header/footer each contain one inline picture:
#   - Pearson Edexcel logo (originally named "image1.png") -> rename to "image2.png"
#   - BTEC logo           (originally named "image2.jpg") -> rename to "image1.jpg"
# Both the "first page" and "default" header/footer instances need the same rename.
#
# InlineShape has no settable Name property in the Word object model (only the
# floating Shape object does), so each inline picture is temporarily converted
# to a Shape, renamed, and converted back to an InlineShape - exactly what a
# real Word automation script does to rename an inline picture.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-Logo($rangeObj) {
    $count = $rangeObj.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $ishp = $rangeObj.InlineShapes.Item($i)
        $shape = $ishp.ConvertToShape()

        if ($shape.Name -eq "image1.png") {
            $shape.Name = "image2.png"
        } elseif ($shape.Name -eq "image2.jpg") {
            $shape.Name = "image1.jpg"
        }

        $shape.ConvertToInlineShape() | Out-Null
    }
}

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
for ($hi = 1; $hi -le 2; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    Rename-Logo $hdr.Range
}

for ($fi = 1; $fi -le 2; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    Rename-Logo $ftr.Range
}
